$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new rows of data following the same pattern as the existing
# "daily" rows at the bottom of the sheet (row 68 = 04-10-2021).
# New row 69: 05-10-2021 ; New row 70: 06-10-2021
# Column A: date label (text), Column B: 50000, Column D: 0

$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "05-10-2021"
$ws.Range("A69").Style = "Normal"
$ws.Range("B69").Value = 50000
$ws.Range("D69").Value = 0

$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "06-10-2021"
$ws.Range("A70").Style = "Normal"
$ws.Range("B70").Value = 50000
$ws.Range("D70").Value = 0
